# Auto-generated-assisted Excel COM-interop script
# Applies the 'river update May 2024' diff to Sheet1:
#  - updates ~104 recalculated statistic cells in rows 6-65
#  - appends 13 new data rows (66-78) for the 2019-2023 period

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Part 1: update existing cell values ----
$ws.Cells.Item(6, 7).Value = 0.0025774334332807
$ws.Cells.Item(7, 7).Value = 0.0025774334332807
$ws.Cells.Item(8, 7).Value = 0.008968473395169201
$ws.Cells.Item(8, 12).Value = 0.00189
$ws.Cells.Item(9, 7).Value = 0.008968473395169201
$ws.Cells.Item(9, 12).Value = 0.00189
$ws.Cells.Item(10, 7).Value = 0.0164255297843871
$ws.Cells.Item(10, 9).Value = 0.0373
$ws.Cells.Item(10, 12).Value = 0.0101
$ws.Cells.Item(10, 13).Value = 0.02403
$ws.Cells.Item(10, 14).Value = 0.03107
$ws.Cells.Item(11, 7).Value = 0.0164255297843871
$ws.Cells.Item(11, 9).Value = 0.0373
$ws.Cells.Item(11, 12).Value = 0.0101
$ws.Cells.Item(11, 13).Value = 0.02403
$ws.Cells.Item(11, 14).Value = 0.03107
$ws.Cells.Item(15, 7).Value = 0.0026924105925782
$ws.Cells.Item(15, 13).Value = 0.00233
$ws.Cells.Item(16, 7).Value = 0.0026924105925782
$ws.Cells.Item(16, 13).Value = 0.00233
$ws.Cells.Item(17, 7).Value = 0.009345447266363801
$ws.Cells.Item(18, 7).Value = 0.009345447266363801
$ws.Cells.Item(19, 7).Value = 0.0176438829819853
$ws.Cells.Item(19, 9).Value = 0.04785
$ws.Cells.Item(19, 14).Value = 0.03757
$ws.Cells.Item(20, 7).Value = 0.0176438829819853
$ws.Cells.Item(20, 9).Value = 0.04785
$ws.Cells.Item(20, 14).Value = 0.03757
$ws.Cells.Item(24, 7).Value = 0.0029878075580572
$ws.Cells.Item(24, 13).Value = 0.00369
$ws.Cells.Item(25, 7).Value = 0.0029878075580572
$ws.Cells.Item(25, 13).Value = 0.00369
$ws.Cells.Item(26, 7).Value = 0.009967113933030399
$ws.Cells.Item(27, 7).Value = 0.009967113933030399
$ws.Cells.Item(28, 6).Value = 0.01216
$ws.Cells.Item(28, 7).Value = 0.0181538829819853
$ws.Cells.Item(28, 9).Value = 0.04485
$ws.Cells.Item(29, 6).Value = 0.01216
$ws.Cells.Item(29, 7).Value = 0.0181538829819853
$ws.Cells.Item(29, 9).Value = 0.04485
$ws.Cells.Item(33, 6).Value = 0.00136
$ws.Cells.Item(33, 7).Value = 0.0032787084902805
$ws.Cells.Item(33, 13).Value = 0.0045
$ws.Cells.Item(33, 14).Value = 0.00591
$ws.Cells.Item(34, 6).Value = 0.00136
$ws.Cells.Item(34, 7).Value = 0.0032787084902805
$ws.Cells.Item(34, 13).Value = 0.0045
$ws.Cells.Item(34, 14).Value = 0.00591
$ws.Cells.Item(35, 6).Value = 0.00523
$ws.Cells.Item(35, 7).Value = 0.0102571133709188
$ws.Cells.Item(36, 6).Value = 0.00523
$ws.Cells.Item(36, 7).Value = 0.0102571133709188
$ws.Cells.Item(37, 6).Value = 0.01216
$ws.Cells.Item(37, 7).Value = 0.018837268376024
$ws.Cells.Item(37, 14).Value = 0.04189
$ws.Cells.Item(38, 6).Value = 0.01216
$ws.Cells.Item(38, 7).Value = 0.018837268376024
$ws.Cells.Item(38, 14).Value = 0.04189
$ws.Cells.Item(42, 6).Value = 0.00159
$ws.Cells.Item(42, 7).Value = 0.0033286377024709
$ws.Cells.Item(42, 13).Value = 0.00428
$ws.Cells.Item(42, 14).Value = 0.0058
$ws.Cells.Item(43, 6).Value = 0.00159
$ws.Cells.Item(43, 7).Value = 0.0033286377024709
$ws.Cells.Item(43, 13).Value = 0.00428
$ws.Cells.Item(43, 14).Value = 0.0058
$ws.Cells.Item(44, 7).Value = 0.011912003396787
$ws.Cells.Item(45, 7).Value = 0.011912003396787
$ws.Cells.Item(46, 6).Value = 0.01337
$ws.Cells.Item(46, 7).Value = 0.0197390594932584
$ws.Cells.Item(47, 6).Value = 0.01337
$ws.Cells.Item(47, 7).Value = 0.0197390594932584
$ws.Cells.Item(51, 7).Value = 0.0035844745601803
$ws.Cells.Item(51, 9).Value = 0.00709
$ws.Cells.Item(51, 13).Value = 0.00423
$ws.Cells.Item(51, 14).Value = 0.00572
$ws.Cells.Item(52, 7).Value = 0.0035844745601803
$ws.Cells.Item(52, 9).Value = 0.00709
$ws.Cells.Item(52, 13).Value = 0.00423
$ws.Cells.Item(52, 14).Value = 0.00572
$ws.Cells.Item(53, 7).Value = 0.0136375620692084
$ws.Cells.Item(54, 7).Value = 0.0136375620692084
$ws.Cells.Item(55, 6).Value = 0.01625
$ws.Cells.Item(55, 7).Value = 0.0208801478630326
$ws.Cells.Item(56, 6).Value = 0.01625
$ws.Cells.Item(56, 7).Value = 0.0208801478630326
$ws.Cells.Item(60, 7).Value = 0.0027661803945041
$ws.Cells.Item(60, 9).Value = 0.00579
$ws.Cells.Item(60, 12).Value = 0.00281
$ws.Cells.Item(60, 13).Value = 0.00434
$ws.Cells.Item(60, 14).Value = 0.00565
$ws.Cells.Item(61, 7).Value = 0.0027661803945041
$ws.Cells.Item(61, 9).Value = 0.00579
$ws.Cells.Item(61, 12).Value = 0.00281
$ws.Cells.Item(61, 13).Value = 0.00434
$ws.Cells.Item(61, 14).Value = 0.00565
$ws.Cells.Item(62, 7).Value = 0.0137416436399872
$ws.Cells.Item(63, 7).Value = 0.0137416436399872
$ws.Cells.Item(64, 6).Value = 0.01625
$ws.Cells.Item(64, 7).Value = 0.0194269831469677
$ws.Cells.Item(64, 13).Value = 0.03345
$ws.Cells.Item(65, 6).Value = 0.01625
$ws.Cells.Item(65, 7).Value = 0.0194269831469677
$ws.Cells.Item(65, 13).Value = 0.03345

# ---- Part 2: append new rows 66-78 (2019-2023 period) ----
# Row 66
$ws.Cells.Item(66, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(66, 2).Value = "Chlorophyll A (92nd Percentile)"
$ws.Cells.Item(66, 3).Value = "A"
$ws.Cells.Item(66, 4).Value = "2019 - 2023"
$ws.Cells.Item(66, 5).Value = "RepSite"
$ws.Cells.Item(66, 6).Value = 3.625
$ws.Cells.Item(66, 7).Value = 5.82592592592593
$ws.Cells.Item(66, 8).Value = 37
$ws.Cells.Item(66, 9).Value = 16.8
$ws.Cells.Item(66, 12).Value = 3.225
$ws.Cells.Item(66, 13).Value = 9.16
$ws.Cells.Item(66, 14).Value = 15.59
$ws.Cells.Item(66, 15).Value = 1825125.6
$ws.Cells.Item(66, 16).Value = 5676951.3
$ws.Cells.Item(66, 17).Value = "Ruapehu District"
$ws.Cells.Item(66, 18).Value = "Whanganui"
$ws.Cells.Item(66, 19).Value = "Upper Whanganui"
$ws.Cells.Item(66, 20).Value = "Whai_1"
$ws.Cells.Item(66, 21).Value = "mg chl-a /m2"

# Row 67
$ws.Cells.Item(67, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(67, 2).Value = "DRP (95th Percentile)"
$ws.Cells.Item(67, 3).Value = "C"
$ws.Cells.Item(67, 4).Value = "2019 - 2023"
$ws.Cells.Item(67, 5).Value = "RepSite"
$ws.Cells.Item(67, 6).Value = 0.03
$ws.Cells.Item(67, 7).Value = 0.0288653846153846
$ws.Cells.Item(67, 8).Value = 0.055
$ws.Cells.Item(67, 9).Value = 0.041
$ws.Cells.Item(67, 12).Value = 0.0315
$ws.Cells.Item(67, 13).Value = 0.03732
$ws.Cells.Item(67, 14).Value = 0.04034
$ws.Cells.Item(67, 15).Value = 1825125.6
$ws.Cells.Item(67, 16).Value = 5676951.3
$ws.Cells.Item(67, 17).Value = "Ruapehu District"
$ws.Cells.Item(67, 18).Value = "Whanganui"
$ws.Cells.Item(67, 19).Value = "Upper Whanganui"
$ws.Cells.Item(67, 20).Value = "Whai_1"
$ws.Cells.Item(67, 21).Value = "mg/L"

# Row 68
$ws.Cells.Item(68, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(68, 2).Value = "DRP (Median)"
$ws.Cells.Item(68, 3).Value = "D"
$ws.Cells.Item(68, 4).Value = "2019 - 2023"
$ws.Cells.Item(68, 5).Value = "RepSite"
$ws.Cells.Item(68, 6).Value = 0.03
$ws.Cells.Item(68, 7).Value = 0.0288653846153846
$ws.Cells.Item(68, 8).Value = 0.055
$ws.Cells.Item(68, 9).Value = 0.041
$ws.Cells.Item(68, 12).Value = 0.0315
$ws.Cells.Item(68, 13).Value = 0.03732
$ws.Cells.Item(68, 14).Value = 0.04034
$ws.Cells.Item(68, 15).Value = 1825125.6
$ws.Cells.Item(68, 16).Value = 5676951.3
$ws.Cells.Item(68, 17).Value = "Ruapehu District"
$ws.Cells.Item(68, 18).Value = "Whanganui"
$ws.Cells.Item(68, 19).Value = "Upper Whanganui"
$ws.Cells.Item(68, 20).Value = "Whai_1"
$ws.Cells.Item(68, 21).Value = "mg/L"

# Row 69
$ws.Cells.Item(69, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(69, 2).Value = "Ammoniacal-N (95th Percentile)"
$ws.Cells.Item(69, 3).Value = "A"
$ws.Cells.Item(69, 4).Value = "2019 - 2023"
$ws.Cells.Item(69, 5).Value = "RepSite"
$ws.Cells.Item(69, 6).Value = 0.00257
$ws.Cells.Item(69, 7).Value = 0.0024603849398505
$ws.Cells.Item(69, 8).Value = 0.0059852671384307
$ws.Cells.Item(69, 9).Value = 0.00566
$ws.Cells.Item(69, 12).Value = 0.00226
$ws.Cells.Item(69, 13).Value = 0.00394
$ws.Cells.Item(69, 14).Value = 0.00535
$ws.Cells.Item(69, 15).Value = 1825125.6
$ws.Cells.Item(69, 16).Value = 5676951.3
$ws.Cells.Item(69, 17).Value = "Ruapehu District"
$ws.Cells.Item(69, 18).Value = "Whanganui"
$ws.Cells.Item(69, 19).Value = "Upper Whanganui"
$ws.Cells.Item(69, 20).Value = "Whai_1"
$ws.Cells.Item(69, 21).Value = "mg NH4-N/L"

# Row 70
$ws.Cells.Item(70, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(70, 2).Value = "Ammoniacal-N (Median)"
$ws.Cells.Item(70, 3).Value = "A"
$ws.Cells.Item(70, 4).Value = "2019 - 2023"
$ws.Cells.Item(70, 5).Value = "RepSite"
$ws.Cells.Item(70, 6).Value = 0.00257
$ws.Cells.Item(70, 7).Value = 0.0024603849398505
$ws.Cells.Item(70, 8).Value = 0.0059852671384307
$ws.Cells.Item(70, 9).Value = 0.00566
$ws.Cells.Item(70, 12).Value = 0.00226
$ws.Cells.Item(70, 13).Value = 0.00394
$ws.Cells.Item(70, 14).Value = 0.00535
$ws.Cells.Item(70, 15).Value = 1825125.6
$ws.Cells.Item(70, 16).Value = 5676951.3
$ws.Cells.Item(70, 17).Value = "Ruapehu District"
$ws.Cells.Item(70, 18).Value = "Whanganui"
$ws.Cells.Item(70, 19).Value = "Upper Whanganui"
$ws.Cells.Item(70, 20).Value = "Whai_1"
$ws.Cells.Item(70, 21).Value = "mg NH4-N/L"

# Row 71
$ws.Cells.Item(71, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(71, 2).Value = "Nitrate-N (95th Percentile)"
$ws.Cells.Item(71, 3).Value = "A"
$ws.Cells.Item(71, 4).Value = "2019 - 2023"
$ws.Cells.Item(71, 5).Value = "RepSite"
$ws.Cells.Item(71, 6).Value = 0.013
$ws.Cells.Item(71, 7).Value = 0.0139724128707565
$ws.Cells.Item(71, 8).Value = 0.056
$ws.Cells.Item(71, 9).Value = 0.0369
$ws.Cells.Item(71, 12).Value = 0.0035
$ws.Cells.Item(71, 13).Value = 0.02462
$ws.Cells.Item(71, 14).Value = 0.03336
$ws.Cells.Item(71, 15).Value = 1825125.6
$ws.Cells.Item(71, 16).Value = 5676951.3
$ws.Cells.Item(71, 17).Value = "Ruapehu District"
$ws.Cells.Item(71, 18).Value = "Whanganui"
$ws.Cells.Item(71, 19).Value = "Upper Whanganui"
$ws.Cells.Item(71, 20).Value = "Whai_1"
$ws.Cells.Item(71, 21).Value = "mg NO3-N/L"

# Row 72
$ws.Cells.Item(72, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(72, 2).Value = "Nitrate-N (Median)"
$ws.Cells.Item(72, 3).Value = "A"
$ws.Cells.Item(72, 4).Value = "2019 - 2023"
$ws.Cells.Item(72, 5).Value = "RepSite"
$ws.Cells.Item(72, 6).Value = 0.013
$ws.Cells.Item(72, 7).Value = 0.0139724128707565
$ws.Cells.Item(72, 8).Value = 0.056
$ws.Cells.Item(72, 9).Value = 0.0369
$ws.Cells.Item(72, 12).Value = 0.0035
$ws.Cells.Item(72, 13).Value = 0.02462
$ws.Cells.Item(72, 14).Value = 0.03336
$ws.Cells.Item(72, 15).Value = 1825125.6
$ws.Cells.Item(72, 16).Value = 5676951.3
$ws.Cells.Item(72, 17).Value = "Ruapehu District"
$ws.Cells.Item(72, 18).Value = "Whanganui"
$ws.Cells.Item(72, 19).Value = "Upper Whanganui"
$ws.Cells.Item(72, 20).Value = "Whai_1"
$ws.Cells.Item(72, 21).Value = "mg NO3-N/L"

# Row 73
$ws.Cells.Item(73, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(73, 2).Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Cells.Item(73, 4).Value = "2019 - 2023"
$ws.Cells.Item(73, 5).Value = "RepSite"
$ws.Cells.Item(73, 6).Value = 0.01625
$ws.Cells.Item(73, 7).Value = 0.0192058293008138
$ws.Cells.Item(73, 8).Value = 0.059
$ws.Cells.Item(73, 9).Value = 0.042
$ws.Cells.Item(73, 12).Value = 0.011
$ws.Cells.Item(73, 13).Value = 0.02799
$ws.Cells.Item(73, 14).Value = 0.04068
$ws.Cells.Item(73, 15).Value = 1825125.6
$ws.Cells.Item(73, 16).Value = 5676951.3
$ws.Cells.Item(73, 17).Value = "Ruapehu District"
$ws.Cells.Item(73, 18).Value = "Whanganui"
$ws.Cells.Item(73, 19).Value = "Upper Whanganui"
$ws.Cells.Item(73, 20).Value = "Whai_1"
$ws.Cells.Item(73, 21).Value = "g/m3"

# Row 74
$ws.Cells.Item(74, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(74, 2).Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Cells.Item(74, 4).Value = "2019 - 2023"
$ws.Cells.Item(74, 5).Value = "RepSite"
$ws.Cells.Item(74, 6).Value = 0.01625
$ws.Cells.Item(74, 7).Value = 0.0192058293008138
$ws.Cells.Item(74, 8).Value = 0.059
$ws.Cells.Item(74, 9).Value = 0.042
$ws.Cells.Item(74, 12).Value = 0.011
$ws.Cells.Item(74, 13).Value = 0.02799
$ws.Cells.Item(74, 14).Value = 0.04068
$ws.Cells.Item(74, 15).Value = 1825125.6
$ws.Cells.Item(74, 16).Value = 5676951.3
$ws.Cells.Item(74, 17).Value = "Ruapehu District"
$ws.Cells.Item(74, 18).Value = "Whanganui"
$ws.Cells.Item(74, 19).Value = "Upper Whanganui"
$ws.Cells.Item(74, 20).Value = "Whai_1"
$ws.Cells.Item(74, 21).Value = "g/m3"

# Row 75
$ws.Cells.Item(75, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(75, 2).Value = "Total Nitrogen (95th Percentile)"
$ws.Cells.Item(75, 4).Value = "2019 - 2023"
$ws.Cells.Item(75, 5).Value = "RepSite"
$ws.Cells.Item(75, 6).Value = 0.04
$ws.Cells.Item(75, 7).Value = 0.0447038470831317
$ws.Cells.Item(75, 8).Value = 0.19
$ws.Cells.Item(75, 9).Value = 0.115
$ws.Cells.Item(75, 12).Value = 0.03
$ws.Cells.Item(75, 13).Value = 0.06
$ws.Cells.Item(75, 14).Value = 0.08599999999999999
$ws.Cells.Item(75, 15).Value = 1825125.6
$ws.Cells.Item(75, 16).Value = 5676951.3
$ws.Cells.Item(75, 17).Value = "Ruapehu District"
$ws.Cells.Item(75, 18).Value = "Whanganui"
$ws.Cells.Item(75, 19).Value = "Upper Whanganui"
$ws.Cells.Item(75, 20).Value = "Whai_1"
$ws.Cells.Item(75, 21).Value = "g/m3"

# Row 76
$ws.Cells.Item(76, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(76, 2).Value = "Total Nitrogen (Median)"
$ws.Cells.Item(76, 4).Value = "2019 - 2023"
$ws.Cells.Item(76, 5).Value = "RepSite"
$ws.Cells.Item(76, 6).Value = 0.04
$ws.Cells.Item(76, 7).Value = 0.0447038470831317
$ws.Cells.Item(76, 8).Value = 0.19
$ws.Cells.Item(76, 9).Value = 0.115
$ws.Cells.Item(76, 12).Value = 0.03
$ws.Cells.Item(76, 13).Value = 0.06
$ws.Cells.Item(76, 14).Value = 0.08599999999999999
$ws.Cells.Item(76, 15).Value = 1825125.6
$ws.Cells.Item(76, 16).Value = 5676951.3
$ws.Cells.Item(76, 17).Value = "Ruapehu District"
$ws.Cells.Item(76, 18).Value = "Whanganui"
$ws.Cells.Item(76, 19).Value = "Upper Whanganui"
$ws.Cells.Item(76, 20).Value = "Whai_1"
$ws.Cells.Item(76, 21).Value = "g/m3"

# Row 77
$ws.Cells.Item(77, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(77, 2).Value = "Total Phosphorus (95th Percentile)"
$ws.Cells.Item(77, 4).Value = "2019 - 2023"
$ws.Cells.Item(77, 5).Value = "RepSite"
$ws.Cells.Item(77, 6).Value = 0.0345
$ws.Cells.Item(77, 7).Value = 0.0352
$ws.Cells.Item(77, 8).Value = 0.07000000000000001
$ws.Cells.Item(77, 9).Value = 0.0585
$ws.Cells.Item(77, 12).Value = 0.0355
$ws.Cells.Item(77, 13).Value = 0.0457
$ws.Cells.Item(77, 14).Value = 0.0535
$ws.Cells.Item(77, 15).Value = 1825125.6
$ws.Cells.Item(77, 16).Value = 5676951.3
$ws.Cells.Item(77, 17).Value = "Ruapehu District"
$ws.Cells.Item(77, 18).Value = "Whanganui"
$ws.Cells.Item(77, 19).Value = "Upper Whanganui"
$ws.Cells.Item(77, 20).Value = "Whai_1"
$ws.Cells.Item(77, 21).Value = "g/m3"

# Row 78
$ws.Cells.Item(78, 1).Value = "Whanganui at Downstream Intake"
$ws.Cells.Item(78, 2).Value = "Total Phosphorus (Median)"
$ws.Cells.Item(78, 4).Value = "2019 - 2023"
$ws.Cells.Item(78, 5).Value = "RepSite"
$ws.Cells.Item(78, 6).Value = 0.0345
$ws.Cells.Item(78, 7).Value = 0.0352
$ws.Cells.Item(78, 8).Value = 0.07000000000000001
$ws.Cells.Item(78, 9).Value = 0.0585
$ws.Cells.Item(78, 12).Value = 0.0355
$ws.Cells.Item(78, 13).Value = 0.0457
$ws.Cells.Item(78, 14).Value = 0.0535
$ws.Cells.Item(78, 15).Value = 1825125.6
$ws.Cells.Item(78, 16).Value = 5676951.3
$ws.Cells.Item(78, 17).Value = "Ruapehu District"
$ws.Cells.Item(78, 18).Value = "Whanganui"
$ws.Cells.Item(78, 19).Value = "Upper Whanganui"
$ws.Cells.Item(78, 20).Value = "Whai_1"
$ws.Cells.Item(78, 21).Value = "g/m3"

Write-Output "Applied river update May 2024: updated recalculated stats + appended rows 66-78"

